# --- Add the "heuristics" worksheet after "multilayer" ---
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("multilayer")
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "heuristics"

# --- Column A: iter # header + numeric iteration counts (real numbers) ---
$ws2.Range("A1").Formula = '="iter #"'
$ws2.Range("A2").Value = 0
$ws2.Range("A3").Value = 10
$ws2.Range("A4").Value = 20
$ws2.Range("A5").Value = 30
$ws2.Range("A6").Value = 40
$ws2.Range("A7").Value = 50
$ws2.Range("A8").Value = 60
$ws2.Range("A9").Value = 70
$ws2.Range("A10").Value = 80
$ws2.Range("A11").Value = 90
$ws2.Range("A12").Value = 100
$ws2.Range("A13").Value = 110
$ws2.Range("A14").Value = 120
$ws2.Range("A15").Value = 130
$ws2.Range("A16").Value = 140
$ws2.Range("A17").Value = 150
$ws2.Range("A18").Value = 160
$ws2.Range("A19").Value = 170
$ws2.Range("A20").Value = 180
$ws2.Range("A21").Value = 190
$ws2.Range("A22").Value = 200
$ws2.Range("A23").Value = 210
$ws2.Range("A24").Value = 220
$ws2.Range("A25").Value = 230
$ws2.Range("A26").Value = 240
$ws2.Range("A27").Value = 250

# --- Columns B, C, G: text content, written via formula then flattened to literal text ---
$ws2.Range("B1").Formula = '="original"'
$ws2.Range("C1").Formula = '="with heuristics"'
$ws2.Range("G1").Formula = '="* ""original"" on every page is the best-performing from previous page. Here is multilayer predictor"'
$ws2.Range("B2").Formula = '="0.5132167547783651"'
$ws2.Range("C2").Formula = '="0.4871899145994307"'
$ws2.Range("G2").Formula = '="Formula extracting last word (printed accuracy number):"'
$ws2.Range("B3").Formula = '="0.535583570557137"'
$ws2.Range("C3").Formula = '="0.5124034160227735"'
$ws2.Range("G3").Formula = '="0 test acc: 0.5132167547783651"'
$ws2.Range("B4").Formula = '="0.5408702724684831"'
$ws2.Range("C4").Formula = '="0.5262301748678324"'
$ws2.Range("G4").Formula = '="10 test acc: 0.535583570557137"'
$ws2.Range("B5").Formula = '="0.5518503456689712"'
$ws2.Range("C5").Formula = '="0.5237901586010574"'
$ws2.Range("G5").Formula = '="20 test acc: 0.5408702724684831"'
$ws2.Range("B6").Formula = '="0.5563237088247255"'
$ws2.Range("C6").Formula = '="0.5685237901586011"'
$ws2.Range("G6").Formula = '="30 test acc: 0.5518503456689712"'
$ws2.Range("B7").Formula = '="0.5673037820252135"'
$ws2.Range("C7").Formula = '="0.555103700691338"'
$ws2.Range("G7").Formula = '="40 test acc: 0.5563237088247255"'
$ws2.Range("B8").Formula = '="0.5734038226921513"'
$ws2.Range("C8").Formula = '="0.590890605937373"'
$ws2.Range("G8").Formula = '="50 test acc: 0.5673037820252135"'
$ws2.Range("B9").Formula = '="0.5860105734038227"'
$ws2.Range("C9").Formula = '="0.6014640097600651"'
$ws2.Range("G9").Formula = '="60 test acc: 0.5734038226921513"'
$ws2.Range("B10").Formula = '="0.588043920292802"'
$ws2.Range("C10").Formula = '="0.6136640910939406"'
$ws2.Range("G10").Formula = '="70 test acc: 0.5860105734038227"'
$ws2.Range("B11").Formula = '="0.6201708011386743"'
$ws2.Range("C11").Formula = '="0.6047173647824319"'
$ws2.Range("G11").Formula = '="80 test acc: 0.588043920292802"'
$ws2.Range("B12").Formula = '="0.6274908499389996"'
$ws2.Range("C12").Formula = '="0.6169174461163074"'
$ws2.Range("G12").Formula = '="90 test acc: 0.6201708011386743"'
$ws2.Range("B13").Formula = '="0.6287108580723871"'
$ws2.Range("C13").Formula = '="0.6295241968279789"'
$ws2.Range("G13").Formula = '="100 test acc: 0.6274908499389996"'
$ws2.Range("B14").Formula = '="0.6384709231394876"'
$ws2.Range("C14").Formula = '="0.6287108580723871"'
$ws2.Range("G14").Formula = '="110 test acc: 0.6287108580723871"'
$ws2.Range("B15").Formula = '="0.6258641724278162"'
$ws2.Range("C15").Formula = '="0.6250508336722245"'
$ws2.Range("G15").Formula = '="120 test acc: 0.6384709231394876"'
$ws2.Range("B16").Formula = '="0.6091907279381863"'
$ws2.Range("C16").Formula = '="0.6352175681171208"'
$ws2.Range("G16").Formula = '="130 test acc: 0.6258641724278162"'
$ws2.Range("B17").Formula = '="0.6425376169174462"'
$ws2.Range("C17").Formula = '="0.6384709231394876"'
$ws2.Range("G17").Formula = '="140 test acc: 0.6091907279381863"'
$ws2.Range("B18").Formula = '="0.6283041886945913"'
$ws2.Range("C18").Formula = '="0.6270841805612037"'
$ws2.Range("G18").Formula = '="150 test acc: 0.6425376169174462"'
$ws2.Range("B19").Formula = '="0.6392842618950793"'
$ws2.Range("C19").Formula = '="0.6360309068727125"'
$ws2.Range("G19").Formula = '="160 test acc: 0.6283041886945913"'
$ws2.Range("B20").Formula = '="0.6457909719398129"'
$ws2.Range("C20").Formula = '="0.6278975193167955"'
$ws2.Range("G20").Formula = '="170 test acc: 0.6392842618950793"'
$ws2.Range("B21").Formula = '="0.6417242781618544"'
$ws2.Range("C21").Formula = '="0.6283041886945913"'
$ws2.Range("G21").Formula = '="180 test acc: 0.6457909719398129"'
$ws2.Range("B22").Formula = '="0.6278975193167955"'
$ws2.Range("C22").Formula = '="0.6278975193167955"'
$ws2.Range("G22").Formula = '="190 test acc: 0.6417242781618544"'
$ws2.Range("B23").Formula = '="0.6421309475396503"'
$ws2.Range("C23").Formula = '="0.6437576250508337"'
$ws2.Range("G23").Formula = '="200 test acc: 0.6278975193167955"'
$ws2.Range("B24").Formula = '="0.6514843432289549"'
$ws2.Range("C24").Formula = '="0.6421309475396503"'
$ws2.Range("G24").Formula = '="210 test acc: 0.6421309475396503"'
$ws2.Range("B25").Formula = '="0.6494509963399756"'
$ws2.Range("C25").Formula = '="0.6445709638064254"'
$ws2.Range("G25").Formula = '="220 test acc: 0.6514843432289549"'
$ws2.Range("B26").Formula = '="0.651077673851159"'
$ws2.Range("C26").Formula = '="0.6437576250508337"'
$ws2.Range("G26").Formula = '="230 test acc: 0.6494509963399756"'
$ws2.Range("B27").Formula = '="0.6364375762505083"'
$ws2.Range("C27").Formula = '="0.6518910126067508"'
$ws2.Range("G27").Formula = '="240 test acc: 0.651077673851159"'
$ws2.Range("G28").Formula = '="250 test acc: 0.6364375762505083"'

# Flatten all the TEXT formulas above into literal string values (keeps them as text, avoids numeric coercion)
$textRng = $ws2.Range("A1:C28,G1:G28")
$textRng.Copy()
$textRng.PasteSpecial(-4163)
$excel.CutCopyMode = $false

# --- Column H: formula extracting the trailing number from column G ---
$ws2.Range("H3").Formula = '=RIGHT(G3, LEN(G3)-FIND("/", SUBSTITUTE(G3," ","/", LEN(G3)-LEN(SUBSTITUTE(G3," ","")))))'
$ws2.Range("H4").Formula = '=RIGHT(G4, LEN(G4)-FIND("/", SUBSTITUTE(G4," ","/", LEN(G4)-LEN(SUBSTITUTE(G4," ","")))))'
$ws2.Range("H5").Formula = '=RIGHT(G5, LEN(G5)-FIND("/", SUBSTITUTE(G5," ","/", LEN(G5)-LEN(SUBSTITUTE(G5," ","")))))'
$ws2.Range("H6").Formula = '=RIGHT(G6, LEN(G6)-FIND("/", SUBSTITUTE(G6," ","/", LEN(G6)-LEN(SUBSTITUTE(G6," ","")))))'
$ws2.Range("H7").Formula = '=RIGHT(G7, LEN(G7)-FIND("/", SUBSTITUTE(G7," ","/", LEN(G7)-LEN(SUBSTITUTE(G7," ","")))))'
$ws2.Range("H8").Formula = '=RIGHT(G8, LEN(G8)-FIND("/", SUBSTITUTE(G8," ","/", LEN(G8)-LEN(SUBSTITUTE(G8," ","")))))'
$ws2.Range("H9").Formula = '=RIGHT(G9, LEN(G9)-FIND("/", SUBSTITUTE(G9," ","/", LEN(G9)-LEN(SUBSTITUTE(G9," ","")))))'
$ws2.Range("H10").Formula = '=RIGHT(G10, LEN(G10)-FIND("/", SUBSTITUTE(G10," ","/", LEN(G10)-LEN(SUBSTITUTE(G10," ","")))))'
$ws2.Range("H11").Formula = '=RIGHT(G11, LEN(G11)-FIND("/", SUBSTITUTE(G11," ","/", LEN(G11)-LEN(SUBSTITUTE(G11," ","")))))'
$ws2.Range("H12").Formula = '=RIGHT(G12, LEN(G12)-FIND("/", SUBSTITUTE(G12," ","/", LEN(G12)-LEN(SUBSTITUTE(G12," ","")))))'
$ws2.Range("H13").Formula = '=RIGHT(G13, LEN(G13)-FIND("/", SUBSTITUTE(G13," ","/", LEN(G13)-LEN(SUBSTITUTE(G13," ","")))))'
$ws2.Range("H14").Formula = '=RIGHT(G14, LEN(G14)-FIND("/", SUBSTITUTE(G14," ","/", LEN(G14)-LEN(SUBSTITUTE(G14," ","")))))'
$ws2.Range("H15").Formula = '=RIGHT(G15, LEN(G15)-FIND("/", SUBSTITUTE(G15," ","/", LEN(G15)-LEN(SUBSTITUTE(G15," ","")))))'
$ws2.Range("H16").Formula = '=RIGHT(G16, LEN(G16)-FIND("/", SUBSTITUTE(G16," ","/", LEN(G16)-LEN(SUBSTITUTE(G16," ","")))))'
$ws2.Range("H17").Formula = '=RIGHT(G17, LEN(G17)-FIND("/", SUBSTITUTE(G17," ","/", LEN(G17)-LEN(SUBSTITUTE(G17," ","")))))'
$ws2.Range("H18").Formula = '=RIGHT(G18, LEN(G18)-FIND("/", SUBSTITUTE(G18," ","/", LEN(G18)-LEN(SUBSTITUTE(G18," ","")))))'
$ws2.Range("H19").Formula = '=RIGHT(G19, LEN(G19)-FIND("/", SUBSTITUTE(G19," ","/", LEN(G19)-LEN(SUBSTITUTE(G19," ","")))))'
$ws2.Range("H20").Formula = '=RIGHT(G20, LEN(G20)-FIND("/", SUBSTITUTE(G20," ","/", LEN(G20)-LEN(SUBSTITUTE(G20," ","")))))'
$ws2.Range("H21").Formula = '=RIGHT(G21, LEN(G21)-FIND("/", SUBSTITUTE(G21," ","/", LEN(G21)-LEN(SUBSTITUTE(G21," ","")))))'
$ws2.Range("H22").Formula = '=RIGHT(G22, LEN(G22)-FIND("/", SUBSTITUTE(G22," ","/", LEN(G22)-LEN(SUBSTITUTE(G22," ","")))))'
$ws2.Range("H23").Formula = '=RIGHT(G23, LEN(G23)-FIND("/", SUBSTITUTE(G23," ","/", LEN(G23)-LEN(SUBSTITUTE(G23," ","")))))'
$ws2.Range("H24").Formula = '=RIGHT(G24, LEN(G24)-FIND("/", SUBSTITUTE(G24," ","/", LEN(G24)-LEN(SUBSTITUTE(G24," ","")))))'
$ws2.Range("H25").Formula = '=RIGHT(G25, LEN(G25)-FIND("/", SUBSTITUTE(G25," ","/", LEN(G25)-LEN(SUBSTITUTE(G25," ","")))))'
$ws2.Range("H26").Formula = '=RIGHT(G26, LEN(G26)-FIND("/", SUBSTITUTE(G26," ","/", LEN(G26)-LEN(SUBSTITUTE(G26," ","")))))'
$ws2.Range("H27").Formula = '=RIGHT(G27, LEN(G27)-FIND("/", SUBSTITUTE(G27," ","/", LEN(G27)-LEN(SUBSTITUTE(G27," ","")))))'
$ws2.Range("H28").Formula = '=RIGHT(G28, LEN(G28)-FIND("/", SUBSTITUTE(G28," ","/", LEN(G28)-LEN(SUBSTITUTE(G28," ","")))))'

# --- Column widths (match content-driven autofit like the original author) ---
$ws2.Columns("B:C").AutoFit()

# --- Selections / active sheet ---
$ws1.Range("E2:F3").Select()
$ws2.Range("C8").Select()
